# Apply the "done up to and including bbr checks" edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1) Update DIR.PRT comment (D4) and bump its row height to fit the wrapped text.
$ws.Range("D4").Value = "Leave empty if PRT are already in BV participant directories. Will search BV directory before searching this directory."
$ws.Rows(4).RowHeight = 30

# 2) Rename the PRT naming convention Field_ID from PRT.FORMAT to PRT.NAMING.
$ws.Range("E7").Value = "PRT.NAMING"

# 3) Insert two new rows (8 and 9) for the Anatomical/Functional naming fields,
#    pushing the existing "Runs" section (and everything after it) down by two.
$ws.Rows("8:9").Insert()

$ws.Range("B8").Value = "Anatomical name"
$ws.Range("C8").Value = "Anat"
$ws.Range("D8").Value = "Name entered in the Create Document Workflow."
$ws.Range("E8").Value = "VMR.NAME"

$ws.Range("B9").Value = "Function name"
$ws.Range("C9").Value = "Func"
$ws.Range("D9").Value = "Name entered in the Create Document Workflow."
$ws.Range("E9").Value = "VTC.NAME"

# 4) Widen column C so the new short codes / existing values are easier to read.
$ws.Columns("C").ColumnWidth = 50.5

# 5) Grow the table / autofilter to cover the two new rows.
$tbl.Resize($ws.Range("A1:E32"))

# 6) Leave the selection where the author left off.
$ws.Range("D4").Select()
